$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3933.3333
$ws.Range("I62").Value = 1700
$ws.Range("J62").Value = 5050
$ws.Range("K62").Value = 1700
$ws.Range("L62").Value = 5050
$ws.Range("M62").Value = -1076
$ws.Range("N62").Value = -6298

$ws.Range("H65").Value = 3933.3333
$ws.Range("I65").Value = 1700
$ws.Range("J65").Value = 5050
$ws.Range("K65").Value = 8500
$ws.Range("L65").Value = 25250
$ws.Range("M65").Value = -5380
$ws.Range("N65").Value = -31490

$ws.Range("H76").Value = 3331
$ws.Range("I76").Value = 3331
$ws.Range("K76").Value = 3331
$ws.Range("M76").Value = -3016

$ws.Range("H79").Value = 3331
$ws.Range("I79").Value = 3331
$ws.Range("K79").Value = 3331
$ws.Range("M79").Value = -2239

$ws.Range("H93").Value = 32331.316
$ws.Range("J93").Value = 32331.316
$ws.Range("L93").Value = 32331.316
$ws.Range("N93").Value = -37323.316

$ws.Range("H132").Value = 405645.4
$ws.Range("I132").Value = 6206.45
$ws.Range("K132").Value = 18619.35
$ws.Range("M132").Value = -16089.35

$ws.Range("H138").Value = 3642.05
$ws.Range("I138").Value = 647.0333000000001
$ws.Range("J138").Value = 4925.6284
$ws.Range("K138").Value = 1941.0999
$ws.Range("L138").Value = 14776.8852
$ws.Range("M138").Value = 3198.9001
$ws.Range("N138").Value = -25056.8852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4742.036
$ws.Range("I32").Value = 4820.15
$ws.Range("J32").Value = 4533.7334
$ws.Range("K32").Value = 4820.15
$ws.Range("L32").Value = 4533.7334
$ws.Range("M32").Value = -4533.15
$ws.Range("N32").Value = -5107.7334

$ws.Range("H61").Value = 1583.258
$ws.Range("I61").Value = 1617.4348
$ws.Range("J61").Value = 1485
$ws.Range("K61").Value = 1617.4348
$ws.Range("L61").Value = 1485
$ws.Range("M61").Value = -1405.4348
$ws.Range("N61").Value = -1909

$ws.Range("H74").Value = 4861.48
$ws.Range("I74").Value = 5900.7334
$ws.Range("K74").Value = 5900.7334
$ws.Range("M74").Value = -5026.7334

$ws.Range("H77").Value = 4861.48
$ws.Range("I77").Value = 5900.7334
$ws.Range("K77").Value = 29503.667
$ws.Range("M77").Value = -25135.667

$ws.Range("H110").Value = 942.4400000000001
$ws.Range("I110").Value = 982.087
$ws.Range("K110").Value = 982.087
$ws.Range("M110").Value = 1062.913

$ws.Range("H136").Value = 1583.258
$ws.Range("I136").Value = 1617.4348
$ws.Range("J136").Value = 1485
$ws.Range("K136").Value = 4852.3044
$ws.Range("L136").Value = 4455
$ws.Range("M136").Value = -2302.3044
$ws.Range("N136").Value = -9555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1237.1111
$ws.Range("I107").Value = 1304.7142
$ws.Range("J107").Value = 1000.5
$ws.Range("K107").Value = 1304.7142
$ws.Range("L107").Value = 1000.5
$ws.Range("M107").Value = 615.2858000000001
$ws.Range("N107").Value = -4840.5

$ws.Range("H132").Value = 50706.668
$ws.Range("J132").Value = 50706.668
$ws.Range("L132").Value = 50706.668
$ws.Range("N132").Value = -60826.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6537170.5
$ws.Range("I16").Value = 12346734
$ws.Range("J16").Value = 1412.125
$ws.Range("K16").Value = 12346734
$ws.Range("L16").Value = 1412.125
$ws.Range("M16").Value = -12346447
$ws.Range("N16").Value = -1986.125

$ws.Range("H31").Value = 12822562
$ws.Range("I31").Value = 789.88464
$ws.Range("K31").Value = 789.88464
$ws.Range("M31").Value = -494.88464

$ws.Range("H34").Value = 12822562
$ws.Range("I34").Value = 789.88464
$ws.Range("K34").Value = 789.88464
$ws.Range("M34").Value = -587.88464

$ws.Range("H58").Value = 1499.5568
$ws.Range("I58").Value = 1461.7222
$ws.Range("J58").Value = 1669.8125
$ws.Range("K58").Value = 1461.7222
$ws.Range("L58").Value = 1669.8125
$ws.Range("M58").Value = -1258.7222
$ws.Range("N58").Value = -2075.8125

$ws.Range("H113").Value = 6537170.5
$ws.Range("I113").Value = 12346734
$ws.Range("J113").Value = 1412.125
$ws.Range("K113").Value = 12346734
$ws.Range("L113").Value = 1412.125
$ws.Range("M113").Value = -12344564
$ws.Range("N113").Value = -5752.125

$ws.Range("H134").Value = 4579.9697
$ws.Range("I134").Value = 6096.6665
$ws.Range("J134").Value = 2759.9333
$ws.Range("K134").Value = 18289.9995
$ws.Range("L134").Value = 8279.7999
$ws.Range("M134").Value = -15754.9995
$ws.Range("N134").Value = -13349.7999

$ws.Range("H136").Value = 1499.5568
$ws.Range("I136").Value = 1461.7222
$ws.Range("J136").Value = 1669.8125
$ws.Range("K136").Value = 4385.1666
$ws.Range("L136").Value = 5009.4375
$ws.Range("M136").Value = -1835.1666
$ws.Range("N136").Value = -10109.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 301.33334
$ws.Range("J97").Value = 301.33334
$ws.Range("L97").Value = 904.0000200000001
$ws.Range("N97").Value = -1896.00002

$ws.Range("H98").Value = 267.33334
$ws.Range("I98").Value = 249
$ws.Range("J98").Value = 304
$ws.Range("K98").Value = 747
$ws.Range("L98").Value = 912
$ws.Range("M98").Value = 751
$ws.Range("N98").Value = -3908

$ws.Range("H107").Value = 62973.875
$ws.Range("I107").Value = 440.91666
$ws.Range("K107").Value = 1322.74998
$ws.Range("M107").Value = 597.2500199999999

$ws.Range("H132").Value = 2636.0715
$ws.Range("J132").Value = 2975.4167
$ws.Range("L132").Value = 26778.7503
$ws.Range("N132").Value = -31838.7503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1862.2307
$ws.Range("I102").Value = 1370.3478
$ws.Range("J102").Value = 5633.3335
$ws.Range("K102").Value = 1370.3478
$ws.Range("L102").Value = 5633.3335
$ws.Range("M102").Value = 251.6522
$ws.Range("N102").Value = -8877.333500000001

$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178

$ws.Range("H122").Value = 2694.6667
$ws.Range("I122").Value = 2079.4
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 6238.200000000001
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -3788.200000000001
$ws.Range("N122").Value = -49900

$ws.Range("H132").Value = 2289.7568
$ws.Range("I132").Value = 1793.6666
$ws.Range("J132").Value = 2759.7368
$ws.Range("K132").Value = 5380.9998
$ws.Range("L132").Value = 8279.2104
$ws.Range("M132").Value = -2850.9998
$ws.Range("N132").Value = -13339.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3939.5173
$ws.Range("I40").Value = 3397.923
$ws.Range("J40").Value = 8633.333000000001
$ws.Range("K40").Value = 3397.923
$ws.Range("L40").Value = 8633.333000000001
$ws.Range("M40").Value = -3261.923
$ws.Range("N40").Value = -8905.333000000001

$ws.Range("H54").Value = 33815.5
$ws.Range("J54").Value = 33815.5
$ws.Range("L54").Value = 33815.5
$ws.Range("N54").Value = -35103.5

$ws.Range("H82").Value = 1507.0869
$ws.Range("I82").Value = 655.25
$ws.Range("J82").Value = 2436.3635
$ws.Range("K82").Value = 655.25
$ws.Range("L82").Value = 2436.3635
$ws.Range("M82").Value = -294.25
$ws.Range("N82").Value = -3158.3635

$ws.Range("H85").Value = 1507.0869
$ws.Range("I85").Value = 655.25
$ws.Range("J85").Value = 2436.3635
$ws.Range("K85").Value = 655.25
$ws.Range("L85").Value = 2436.3635
$ws.Range("M85").Value = 592.75
$ws.Range("N85").Value = -4932.363499999999

$ws.Range("H122").Value = 2861.611
$ws.Range("I122").Value = 1675.3334
$ws.Range("K122").Value = 5026.0002
$ws.Range("M122").Value = -2576.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 33366.668
$ws.Range("J92").Value = 33366.668
$ws.Range("L92").Value = 33366.668
$ws.Range("N92").Value = -38358.668
